$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.91
$ws.Range("H2").Value = 3.25
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("U2").Value = 3.95
$ws.Range("W2").Value = 4.5
$ws.Range("X2").Value = 1.18
$ws.Range("AA2").Value = 2.1
$ws.Range("AB2").Value = 1.67
$ws.Range("AD2").Value = 8
$ws.Range("AK2").Value = 19
# Row 3
$ws.Range("G3").Value = 2.7
$ws.Range("I3").Value = 3.6
$ws.Range("J3").Value = 3.75
$ws.Range("K3").Value = 1.67
$ws.Range("N3").Value = 3.6
$ws.Range("Q3").Value = 3.4
$ws.Range("R3").Value = 1.32
$ws.Range("U3").Value = 9
$ws.Range("V3").Value = 1.06
$ws.Range("AD3").Value = 10
$ws.Range("AF3").Value = 29
$ws.Range("AL3").Value = 151
$ws.Range("AM3").Value = 6
# Row 7
$ws.Range("G7").Value = 1.75
$ws.Range("H7").Value = 3.4
$ws.Range("I7").Value = 5
$ws.Range("J7").Value = 2.5
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 5.5
$ws.Range("N7").Value = 8
$ws.Range("O7").Value = 1.44
$ws.Range("P7").Value = 2.63
$ws.Range("U7").Value = 3.65
$ws.Range("AF7").Value = 13
$ws.Range("AI7").Value = 7
$ws.Range("AJ7").Value = 6.5
$ws.Range("AN7").Value = 23
$ws.Range("AO7").Value = 17
$ws.Range("AQ7").Value = 41
# Row 8
$ws.Range("G8").Value = 1.91
$ws.Range("M8").Value = 1.1
$ws.Range("N8").Value = 7
$ws.Range("O8").Value = 1.44
$ws.Range("P8").Value = 2.63
$ws.Range("Q8").Value = 1.8
$ws.Range("R8").Value = 2.05
$ws.Range("S8").Value = 2.4
$ws.Range("T8").Value = 1.53
$ws.Range("U8").Value = 3.6
$ws.Range("V8").Value = 1.28
$ws.Range("W8").Value = 4.5
$ws.Range("X8").Value = 1.18
$ws.Range("AE8").Value = 9.5
$ws.Range("AF8").Value = 17
$ws.Range("AK8").Value = 17
$ws.Range("AM8").Value = 9.5
$ws.Range("AO8").Value = 15
$ws.Range("AP8").Value = 41
# Row 9
$ws.Range("G9").Value = 3.9
$ws.Range("H9").Value = 3.1
$ws.Range("I9").Value = 2.05
$ws.Range("J9").Value = 4.75
$ws.Range("L9").Value = 2.88
$ws.Range("M9").Value = 1.11
$ws.Range("N9").Value = 6.5
$ws.Range("Q9").Value = 1.9
$ws.Range("R9").Value = 1.95
$ws.Range("V9").Value = 1.23
$ws.Range("AD9").Value = 19
$ws.Range("AE9").Value = 15
$ws.Range("AG9").Value = 41
$ws.Range("AH9").Value = 51
$ws.Range("AN9").Value = 8.5
$ws.Range("AO9").Value = 9.5
# Row 10
$ws.Range("G10").Value = 3.25
$ws.Range("H10").Value = 2.45
$ws.Range("I10").Value = 2.75
$ws.Range("J10").Value = 4
$ws.Range("L10").Value = 3.6
$ws.Range("M10").Value = 1.2
$ws.Range("N10").Value = 4.05
$ws.Range("O10").Value = 1.78
$ws.Range("P10").Value = 1.93
$ws.Range("S10").Value = 3.25
$ws.Range("T10").Value = 1.29
$ws.Range("W10").Value = 6.1
$ws.Range("X10").Value = 1.09
$ws.Range("Z10").Value = 2
$ws.Range("AA10").Value = 2.4
$ws.Range("AB10").Value = 1.5
$ws.Range("AC10").Value = 6.1
$ws.Range("AD10").Value = 14.5
$ws.Range("AE10").Value = 13
$ws.Range("AF10").Value = 50
$ws.Range("AG10").Value = 45
$ws.Range("AI10").Value = 4.05
$ws.Range("AK10").Value = 22
$ws.Range("AL10").Value = 175
$ws.Range("AM10").Value = 5.3
$ws.Range("AN10").Value = 11.5
$ws.Range("AO10").Value = 12
$ws.Range("AQ10").Value = 37
$ws.Range("AR10").Value = 70
# Row 11
$ws.Range("G11").Value = 3.2
$ws.Range("I11").Value = 2.65
$ws.Range("J11").Value = 4.15
$ws.Range("L11").Value = 3.45
$ws.Range("P11").Value = 2
$ws.Range("AA11").Value = 2.35
$ws.Range("AC11").Value = 6.2
$ws.Range("AD11").Value = 14.5
$ws.Range("AE11").Value = 13
$ws.Range("AF11").Value = 50
$ws.Range("AG11").Value = 45
$ws.Range("AH11").Value = 70
$ws.Range("AM11").Value = 5.5
$ws.Range("AN11").Value = 11.25
$ws.Range("AO11").Value = 11
$ws.Range("AP11").Value = 32
$ws.Range("AQ11").Value = 32
$ws.Range("AR11").Value = 60
# Row 12
$ws.Range("G12").Value = 3.1
$ws.Range("H12").Value = 3.6
$ws.Range("I12").Value = 2.1
$ws.Range("J12").Value = 3.5
$ws.Range("L12").Value = 2.63
$ws.Range("AC12").Value = 15
$ws.Range("AD12").Value = 19
$ws.Range("AE12").Value = 12
$ws.Range("AF12").Value = 34
$ws.Range("AG12").Value = 23
$ws.Range("AH12").Value = 26
$ws.Range("AJ12").Value = 7.5
$ws.Range("AK12").Value = 12
$ws.Range("AN12").Value = 12
$ws.Range("AO12").Value = 9
$ws.Range("AP12").Value = 21
$ws.Range("AQ12").Value = 15
# Row 16
$ws.Range("G16").Value = 2.45
$ws.Range("I16").Value = 2.7
$ws.Range("J16").Value = 3
$ws.Range("L16").Value = 3.25
$ws.Range("M16").Value = 1.02
$ws.Range("N16").Value = 11
$ws.Range("Y16").Value = 1.4
$ws.Range("AC16").Value = 9
$ws.Range("AE16").Value = 10
$ws.Range("AF16").Value = 23
$ws.Range("AM16").Value = 9.5
$ws.Range("AO16").Value = 11
$ws.Range("AP16").Value = 26
$ws.Range("AS16").Value = 500
# Row 17
$ws.Range("G17").Value = 3.25
$ws.Range("H17").Value = 3.6
$ws.Range("I17").Value = 2
$ws.Range("J17").Value = 3.6
$ws.Range("L17").Value = 2.6
$ws.Range("S17").Value = 1.67
$ws.Range("T17").Value = 2.15
$ws.Range("Y17").Value = 1.33
$ws.Range("AC17").Value = 13
$ws.Range("AD17").Value = 19
$ws.Range("AN17").Value = 11
$ws.Range("AO17").Value = 9
$ws.Range("AP17").Value = 19
$ws.Range("AQ17").Value = 15
$ws.Range("AS17").Value = 151
# Row 18
$ws.Range("I18").Value = 2.5
$ws.Range("K18").Value = 2.25
$ws.Range("L18").Value = 3
$ws.Range("M18").Value = 1.04
$ws.Range("N18").Value = 9
$ws.Range("O18").Value = 1.22
$ws.Range("P18").Value = 4
$ws.Range("S18").Value = 1.7
$ws.Range("T18").Value = 2.1
$ws.Range("U18").Value = 2
$ws.Range("V18").Value = 1.8
$ws.Range("W18").Value = 2.63
$ws.Range("X18").Value = 1.44
$ws.Range("Y18").Value = 1.33
$ws.Range("Z18").Value = 3.25
$ws.Range("AM18").Value = 11
$ws.Range("AS18").Value = 126
# Row 20
$ws.Range("G20").Value = 2.9
$ws.Range("I20").Value = 2.38
$ws.Range("J20").Value = 3.6
$ws.Range("K20").Value = 2.05
$ws.Range("L20").Value = 3.1
$ws.Range("M20").Value = 1.07
$ws.Range("N20").Value = 9
$ws.Range("S20").Value = 2.08
$ws.Range("T20").Value = 1.73
$ws.Range("W20").Value = 3.75
$ws.Range("X20").Value = 1.25
$ws.Range("Y20").Value = 1.44
$ws.Range("Z20").Value = 2.63
$ws.Range("AA20").Value = 1.83
$ws.Range("AB20").Value = 1.83
$ws.Range("AC20").Value = 9
$ws.Range("AD20").Value = 15
$ws.Range("AM20").Value = 7.5
$ws.Range("AN20").Value = 11
# Row 21
$ws.Range("Y21").Value = 1.4
# Row 22
$ws.Range("Y22").Value = 1.36
# Row 23
$ws.Range("G23").Value = 1.78
$ws.Range("H23").Value = 3.75
$ws.Range("I23").Value = 3.95
$ws.Range("J23").Value = 2.3
$ws.Range("K23").Value = 2.32
$ws.Range("L23").Value = 4.1
$ws.Range("W23").Value = 2.32
$ws.Range("Y23").Value = 1.29
$ws.Range("Z23").Value = 3.25
$ws.Range("AA23").Value = 1.53
$ws.Range("AB23").Value = 2.35
$ws.Range("AC23").Value = 9.75
$ws.Range("AD23").Value = 10.25
$ws.Range("AF23").Value = 16
$ws.Range("AG23").Value = 12.5
$ws.Range("AH23").Value = 19
$ws.Range("AJ23").Value = 7.7
$ws.Range("AK23").Value = 12
$ws.Range("AM23").Value = 16
$ws.Range("AN23").Value = 27
$ws.Range("AO23").Value = 13
$ws.Range("AP23").Value = 60
$ws.Range("AQ23").Value = 30
$ws.Range("AS23").Value = 250
# Row 24
$ws.Range("G24").Value = 1.52
$ws.Range("H24").Value = 4.45
$ws.Range("J24").Value = 1.98
$ws.Range("K24").Value = 2.52
$ws.Range("M24").Value = 1.02
$ws.Range("N24").Value = 10
$ws.Range("O24").Value = 1.13
$ws.Range("P24").Value = 5.2
$ws.Range("S24").Value = 1.4
$ws.Range("T24").Value = 2.72
$ws.Range("W24").Value = 1.98
$ws.Range("X24").Value = 1.75
$ws.Range("Y24").Value = 1.25
$ws.Range("Z24").Value = 3.55
$ws.Range("AA24").Value = 1.5
$ws.Range("AB24").Value = 2.42
$ws.Range("AC24").Value = 11.5
$ws.Range("AD24").Value = 10.25
$ws.Range("AE24").Value = 8.5
$ws.Range("AF24").Value = 13
$ws.Range("AG24").Value = 10.75
$ws.Range("AH24").Value = 17
$ws.Range("AI24").Value = 10
$ws.Range("AJ24").Value = 9.5
$ws.Range("AL24").Value = 40
$ws.Range("AM24").Value = 22
$ws.Range("AO24").Value = 16.5
$ws.Range("AS24").Value = 200
